$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44630
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15500
$ws.Range("P2").Value = 861
# Row 3
$ws.Range("D3").Value = 44775
$ws.Range("J3").Value = 100
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17500
$ws.Range("P3").Value = 972
# Row 4
$ws.Range("D4").Value = 44769
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17500
$ws.Range("P4").Value = 972
# Row 5
$ws.Range("D5").Value = 44656
# Row 6
$ws.Range("D6").Value = 44782
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("P6").Value = 972
# Row 7
$ws.Range("D7").Value = 44628
# Row 8
$ws.Range("D8").Value = 44799
$ws.Range("J8").Value = 60
# Row 9
$ws.Range("D9").Value = 44797
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 16000
$ws.Range("L9").Value = 17000
$ws.Range("M9").Value = 16500
$ws.Range("P9").Value = 917
# Row 10
$ws.Range("D10").Value = 44804
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15500
$ws.Range("P10").Value = 861
# Row 11
$ws.Range("D11").Value = 44790
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("P11").Value = 972
# Row 12
$ws.Range("D12").Value = 44791
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17500
$ws.Range("P12").Value = 972
# Row 13
$ws.Range("D13").Value = 44819
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 15000
$ws.Range("P13").Value = 833
# Row 14
$ws.Range("D14").Value = 44818
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 15000
$ws.Range("P14").Value = 833
# Row 15
$ws.Range("D15").Value = 44809
$ws.Range("K15").Value = 14000
$ws.Range("M15").Value = 14500
$ws.Range("P15").Value = 806
# Row 16
$ws.Range("D16").Value = 44811
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("P16").Value = 806
# Row 17
$ws.Range("D17").Value = 44813
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("P17").Value = 806
# Row 18
$ws.Range("D18").Value = 44645
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 15500
$ws.Range("P18").Value = 861
# Row 19
$ws.Range("D19").Value = 44651
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15500
$ws.Range("P19").Value = 861
# Row 20
$ws.Range("D20").Value = 44659
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 80
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15500
$ws.Range("P20").Value = 861
# Row 21
$ws.Range("D21").Value = 44832
$ws.Range("J21").Value = 60
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = 17000
$ws.Range("P21").Value = 944
# Row 22
$ws.Range("D22").Value = 44785
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17500
$ws.Range("P22").Value = 972
# Row 23
$ws.Range("D23").Value = 44637
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 16000
$ws.Range("M23").Value = 15500
$ws.Range("P23").Value = 861
# Row 24
$ws.Range("D24").Value = 44664
$ws.Range("J24").Value = 160
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 16000
$ws.Range("M24").Value = 15500
$ws.Range("P24").Value = 861
# Row 25
$ws.Range("D25").Value = 44839
$ws.Range("J25").Value = 100
# Row 26
$ws.Range("D26").Value = 44847
$ws.Range("J26").Value = 120
$ws.Range("K26").Value = 17000
$ws.Range("M26").Value = 17000
$ws.Range("P26").Value = 944
# Row 27
$ws.Range("D27").Value = 44761
# Row 28
$ws.Range("D28").Value = 44771
$ws.Range("J28").Value = 60
$ws.Range("L28").Value = 18000
$ws.Range("M28").Value = 17500
$ws.Range("P28").Value = 972
# Row 29
$ws.Range("D29").Value = 44635
$ws.Range("J29").Value = 100
# Row 30
$ws.Range("D30").Value = 44642
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 16000
$ws.Range("M30").Value = 15500
$ws.Range("P30").Value = 861
# Row 31
$ws.Range("D31").Value = 44830
$ws.Range("J31").Value = 60
$ws.Range("L31").Value = 17000
$ws.Range("M31").Value = 17000
$ws.Range("P31").Value = 944
# Row 32
$ws.Range("D32").Value = 44763
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 17000
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = 17500
$ws.Range("P32").Value = 972
# Row 33
$ws.Range("D33").Value = 44384
$ws.Range("J33").Value = 120
# Row 34
$ws.Range("D34").Value = 44384
$ws.Range("I34").Value = 'Segunda'
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 15000
$ws.Range("P34").Value = 833
# Row 35
$ws.Range("D35").Value = 44658
$ws.Range("J35").Value = 80
# Row 36
$ws.Range("D36").Value = 44649
$ws.Range("J36").Value = 60
$ws.Range("K36").Value = 15000
$ws.Range("L36").Value = 16000
$ws.Range("M36").Value = 15500
$ws.Range("P36").Value = 861
